# Weekly update of "Fruto del paraíso" hortaliza prices at Macroferia Regional de Talca.
# The data rows (2, 3 and 5) are shifted: the newest record moves into row 2,
# the former row 2 record moves into row 3, and the former row 3 record moves
# into row 5. Row 4 is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: becomes the (previously) newest record (old row 5) ---
$ws.Range("D2").Value = 44291
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 11000
$ws.Range("L2").Value = 11000
$ws.Range("M2").Value = 11000
$ws.Range("P2").Value = 550

# --- Row 3: becomes the (previously) row 2 record ---
$ws.Range("D3").Value = 44284
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 10000
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10000
$ws.Range("P3").Value = 500

# --- Row 5: becomes the (previously) row 3 record ---
$ws.Range("D5").Value = 44277
$ws.Range("J5").Value = 150
